# Finish admin side dashboard: add new orders (rows 7-10) and refresh the
# totals row (now row 12) with updated aggregate figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force a value that Excel would otherwise auto-convert (dates, plain
    # numbers, "$123" amounts, ...) to be stored as literal text, then drop
    # the temporary text number-format so the cell keeps the default style.
    $r = $ws.Range($range)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

# --- Back-fill CouponAmount (column G) with 0 for the two rows that didn't
#     have a discount previously recorded, matching the rest of the table. ---
$ws.Range("G3").Value = 0
$ws.Range("G4").Value = 0

# --- New order: 2024-05-27, wildcraft laptop bag, paid via Wallet ---
$ws.Range("A7").Value = 513307
Set-TextValue "B7" "2024-05-27"
$ws.Range("C7").Value = "wildcraft"
Set-TextValue "D7" "1"
Set-TextValue "E7" "50"
Set-TextValue "F7" "150"
Set-TextValue "H7" "`$150"
$ws.Range("I7").Value = "Wallet"
$ws.Range("J7").Value = "Delivered"

# --- Row 8 used to hold the totals; it now becomes a regular order row for
#     2024-05-27, paid via Online Payment. ---
$ws.Range("A8").Value = 474100
Set-TextValue "B8" "2024-05-27"
$ws.Range("C8").Value = " laptop bag2, Hp laptop bag"
$ws.Range("D8").Value = "3, 5"
$ws.Range("E8").Value = "%, 20"
$ws.Range("F8").Value = "144$, 80"
$ws.Range("G8").Value = 25
Set-TextValue "H8" "`$807"
$ws.Range("I8").Value = "Online Payment"
$ws.Range("J8").Value = "Delivered"

# --- New order: 2024-05-27,  laptop bag3, paid via COD ---
$ws.Range("A9").Value = 639057
Set-TextValue "B9" "2024-05-27"
$ws.Range("C9").Value = " laptop bag3"
Set-TextValue "D9" "5"
Set-TextValue "E9" "0"
Set-TextValue "F9" "200"
$ws.Range("G9").Value = 50
Set-TextValue "H9" "`$950"
$ws.Range("I9").Value = "COD"
$ws.Range("J9").Value = "Delivered"

# --- New order: 2024-05-30, wildcraft, paid via Online Payment ---
$ws.Range("A10").Value = 916703
Set-TextValue "B10" "2024-05-30"
$ws.Range("C10").Value = "wildcraft"
Set-TextValue "D10" "5"
Set-TextValue "F10" "250"
Set-TextValue "H10" "`$1250"
$ws.Range("I10").Value = "Online Payment"
$ws.Range("J10").Value = "Delivered"

# --- Totals row moves from row 8 down to row 12 (row 11 stays blank), with
#     refreshed aggregate figures. ---
$ws.Range("D12").Value = "Total Orders: 9"
$ws.Range("G12").Value = "TotalDiscount: `$175"
$ws.Range("H12").Value = "Total Sales: `$6059"
